# Split the "Variables" sheet's single "position" column into two columns:
# "pivot" (figures / stub / heading) and "order" (the numeric position within
# that pivot group). Close #124.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Variables")

# Make room for the new "pivot" column in front of the old "position"
# column (which becomes "order").
$ws.Columns("A:A").Insert()

# New column A: pivot ("figures" / "stub" / "stub" / "heading")
$ws.Range("A1").Value = "pivot"
$ws.Range("A2").Value = "figures"
$ws.Range("A3").Value = "stub"
$ws.Range("A4").Value = "stub"
$ws.Range("A5").Value = "heading"

# Column B (previously A, "position") is renamed "order" and holds the
# numeric ordering within each pivot group instead of the old h1/s1/s2 codes.
$ws.Range("B1").Value = "order"
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 1

# --- Restore/update the view state (active sheet + selections) recorded in
# the workbook, moving the active tab from "Data" to "Variables". ---
$wsTable = $wb.Worksheets.Item("Table")
$wsTable.Activate()
$wsTable.Range("A89").Select()

$wsCodelists = $wb.Worksheets.Item("Codelists")
$wsCodelists.Activate()
$wsCodelists.Range("M56").Select()

$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("D1").Select()

$ws.Activate()
$ws.Range("F11").Select()
